# Update column F (dSF) values per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = -3
    3  = -2
    4  = 0
    6  = -5
    7  = -3
    8  = -5
    9  = 7
    10 = -3
    11 = -2
    12 = 0
    13 = 2
    14 = 4
    15 = 2
    16 = 6
    17 = 1
    18 = 2
    19 = 7
    20 = -4
    21 = -1
    22 = 0
    23 = 1
    24 = -1
    26 = -1
    27 = 1
    28 = 4
    29 = 4
    30 = 2
    31 = -1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("F$row").Value = $newValues[$row]
}
